# C5-PowerPoint.pptx edit: 2020-05-05 commit
#
# 1) Slide 6's table switches to a different table style GUID.
# 2) The presentation's theme colour palette is swapped from the "Integral"
#    palette to the classic "Office" palette (the Integral palette moves to
#    the notes-master theme part in the canonical OOXML; that part is not
#    reachable through the PowerPoint object model, so here we apply the
#    "Office" colours to the one Theme object the OM exposes - the design
#    used by the slide master/notes master/slides alike).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 ---------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{59C7D487-47D1-4E88-A3F0-7D6A5CDFD4E1}")
    }
}

# --- 2) Theme colours: Integral -> Office ----------------------------------
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

# index -> (name, RGB as R + G*256 + B*65536, matching the classic Office theme)
$officeColors = @{
    1  = 0          # dk1      000000
    2  = 16777215    # lt1      FFFFFF
    3  = 6968388      # dk2      44546A
    4  = 15132391     # lt2      E7E6E6
    5  = 13998939     # accent1  5B9BD5
    6  = 3243501      # accent2  ED7D31
    7  = 10855845     # accent3  A5A5A5
    8  = 49407        # accent4  FFC000
    9  = 12874308     # accent5  4472C4
    10 = 4697456      # accent6  70AD47
    11 = 12673797     # hlink    0563C1
    12 = 7491477      # folHlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $scheme.Colors($i).RGB = $officeColors[$i]
}

$p.Save()
